# The source sheet gains a 3-row x 8-column table of worker/timesheet data.
# Every value (including the ones that look numeric/date-like, e.g. "143",
# "1" and "2023-06-05") must land as plain TEXT - exactly what the target
# workbook stores (shared-string cells, no numeric coercion). Typing a
# numeric-looking string straight into a cell makes Excel auto-detect it as
# a number/date, so each cell is briefly forced to Text format, written,
# and then has that temporary formatting cleared again so no stray
# per-cell style sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Андрей", "Степченков", "143", "Рабочий", "1", "2023-06-05", "16:58", "16:59")

for ($row = 1; $row -le 3; $row++) {
    for ($col = 1; $col -le 8; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $headers[$col - 1]
        $cell.ClearFormats()
    }
}

$ws.Range("A1:J8").Select()
